$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 750000400
$ws.Range("J86").Value = 500000740
$ws.Range("L86").Value = 500000740
$ws.Range("N86").Value = -500002986

$ws.Range("H89").Value = 750000400
$ws.Range("J89").Value = 500000740
$ws.Range("L89").Value = 2500003700
$ws.Range("N89").Value = -2500014932

$ws.Range("H100").Value = 2560.5715
$ws.Range("I100").Value = 1769.4286
$ws.Range("K100").Value = 1769.4286
$ws.Range("M100").Value = -1228.4286

$ws.Range("H131").Value = 2551.3
$ws.Range("I131").Value = 1755.4667
$ws.Range("K131").Value = 5266.4001
$ws.Range("M131").Value = -226.4000999999998

$ws.Range("H132").Value = 1684.5238
$ws.Range("I132").Value = 1519
$ws.Range("K132").Value = 4557
$ws.Range("M132").Value = -2027

$ws.Range("H137").Value = 4339602.5
$ws.Range("I137").Value = 138998.17
$ws.Range("K137").Value = 416994.51
$ws.Range("M137").Value = -414444.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 13336800
$ws.Range("I10").Value = 13336800
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 13336800
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -13336630
$ws.Range("N10").ClearContents()

$ws.Range("H24").Value = 80346.336
$ws.Range("J24").Value = 80346.336
$ws.Range("L24").Value = 80346.336
$ws.Range("N24").Value = -81094.336

$ws.Range("H32").Value = 18519558
$ws.Range("I32").Value = 20409010
$ws.Range("K32").Value = 20409010
$ws.Range("M32").Value = -20408723

$ws.Range("H45").Value = 4219.5527
$ws.Range("I45").Value = 4776.3335
$ws.Range("J45").Value = 3962.577
$ws.Range("K45").Value = 4776.3335
$ws.Range("L45").Value = 3962.577
$ws.Range("M45").Value = -4399.3335
$ws.Range("N45").Value = -4716.577

$ws.Range("H74").Value = 2326.3462
$ws.Range("I74").Value = 2557.2354
$ws.Range("J74").Value = 1890.2222
$ws.Range("K74").Value = 2557.2354
$ws.Range("L74").Value = 1890.2222
$ws.Range("M74").Value = -1683.2354
$ws.Range("N74").Value = -3638.2222

$ws.Range("H77").Value = 2326.3462
$ws.Range("I77").Value = 2557.2354
$ws.Range("J77").Value = 1890.2222
$ws.Range("K77").Value = 12786.177
$ws.Range("L77").Value = 9451.110999999999
$ws.Range("M77").Value = -8418.177
$ws.Range("N77").Value = -18187.111

$ws.Range("H100").Value = 80346.336
$ws.Range("J100").Value = 80346.336
$ws.Range("L100").Value = 80346.336
$ws.Range("N100").Value = -82510.336

$ws.Range("H101").Value = 112597.4
$ws.Range("J101").Value = 112597.4
$ws.Range("L101").Value = 112597.4
$ws.Range("N101").Value = -119087.4

$ws.Range("H104").Value = 96995
$ws.Range("J104").Value = 96995
$ws.Range("L104").Value = 96995
$ws.Range("N104").Value = -103983

$ws.Range("H132").Value = 2502.4866
$ws.Range("I132").Value = 2721.0322
$ws.Range("J132").Value = 1373.3334
$ws.Range("K132").Value = 8163.096600000001
$ws.Range("L132").Value = 4120.0002
$ws.Range("M132").Value = -5633.096600000001
$ws.Range("N132").Value = -9180.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4332.3335
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 4997
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4997
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -5491

$ws.Range("H86").Value = 1280.125
$ws.Range("I86").Value = 1004
$ws.Range("J86").Value = 1740.3334
$ws.Range("K86").Value = 1004
$ws.Range("L86").Value = 1740.3334
$ws.Range("M86").Value = 119
$ws.Range("N86").Value = -3986.3334

$ws.Range("H89").Value = 1280.125
$ws.Range("I89").Value = 1004
$ws.Range("J89").Value = 1740.3334
$ws.Range("K89").Value = 5020
$ws.Range("L89").Value = 8701.666999999999
$ws.Range("M89").Value = 596
$ws.Range("N89").Value = -19933.667

$ws.Range("H99").Value = 2690.9
$ws.Range("I99").Value = 1869.6666
$ws.Range("J99").Value = 3042.8572
$ws.Range("K99").Value = 1869.6666
$ws.Range("L99").Value = 3042.8572
$ws.Range("M99").Value = -371.6666
$ws.Range("N99").Value = -6038.8572

$ws.Range("H107").Value = 3927
$ws.Range("I107").Value = 3908.6155
$ws.Range("J107").Value = 3986.75
$ws.Range("K107").Value = 3908.6155
$ws.Range("L107").Value = 3986.75
$ws.Range("M107").Value = -1988.6155
$ws.Range("N107").Value = -7826.75

$ws.Range("H134").Value = 3816.3333
$ws.Range("I134").Value = 3731.7856
$ws.Range("K134").Value = 11195.3568
$ws.Range("M134").Value = -8660.356800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7714.2856
$ws.Range("J22").Value = 7499.5
$ws.Range("L22").Value = 7499.5
$ws.Range("N22").Value = -8199.5

$ws.Range("H31").Value = 4641.2954
$ws.Range("J31").Value = 5145.543
$ws.Range("L31").Value = 5145.543
$ws.Range("N31").Value = -5735.543

$ws.Range("H34").Value = 4641.2954
$ws.Range("J34").Value = 5145.543
$ws.Range("L34").Value = 5145.543
$ws.Range("N34").Value = -5549.543

$ws.Range("H122").Value = 3868.6667
$ws.Range("I122").Value = 2490.6428
$ws.Range("K122").Value = 7471.928400000001
$ws.Range("M122").Value = -5021.928400000001

$ws.Range("H132").Value = 3254
$ws.Range("I132").Value = 4306
$ws.Range("J132").Value = 1150
$ws.Range("K132").Value = 12918
$ws.Range("L132").Value = 3450
$ws.Range("M132").Value = -10388
$ws.Range("N132").Value = -8510

$ws.Range("H134").Value = 863.3333
$ws.Range("I134").Value = 863.3333
$ws.Range("K134").Value = 2589.9999
$ws.Range("M134").Value = -54.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 500
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -1230
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 500
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -564
$ws.Range("N67").ClearContents()

$ws.Range("H113").Value = 1803.4783
$ws.Range("I113").Value = 886.125
$ws.Range("J113").Value = 2292.7334
$ws.Range("K113").Value = 2658.375
$ws.Range("L113").Value = 6878.2002
$ws.Range("M113").Value = -488.375
$ws.Range("N113").Value = -11218.2002

$ws.Range("H131").Value = 1586.0577
$ws.Range("J131").Value = 1687
$ws.Range("L131").Value = 5061
$ws.Range("N131").Value = -15141

$ws.Range("H140").Value = 2057
$ws.Range("I140").Value = 2057
$ws.Range("K140").Value = 6171
$ws.Range("M140").Value = -991

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 6699998.5

$ws.Range("H95").Value = 64950
$ws.Range("J95").Value = 64950
$ws.Range("L95").Value = 64950
$ws.Range("N95").Value = -70442

$ws.Range("H102").Value = 1531.5358
$ws.Range("I102").Value = 1424.826
$ws.Range("J102").Value = 2022.4
$ws.Range("K102").Value = 1424.826
$ws.Range("L102").Value = 2022.4
$ws.Range("M102").Value = 197.174
$ws.Range("N102").Value = -5266.4

$ws.Range("H107").Value = 1267.5
$ws.Range("I107").Value = 1095.1818
$ws.Range("K107").Value = 1095.1818
$ws.Range("M107").Value = 824.8181999999999

$ws.Range("H126").Value = 2471.8262
$ws.Range("I126").Value = 1851.7858
$ws.Range("K126").Value = 5555.357400000001
$ws.Range("M126").Value = -3085.357400000001

$ws.Range("H132").Value = 3348.2068
$ws.Range("I132").Value = 3057.842
$ws.Range("J132").Value = 3899.9
$ws.Range("K132").Value = 9173.526
$ws.Range("L132").Value = 11699.7
$ws.Range("M132").Value = -6643.526
$ws.Range("N132").Value = -16759.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2358.0908
$ws.Range("I7").Value = 2351.4443
$ws.Range("K7").Value = 2351.4443
$ws.Range("M7").Value = -2239.4443

$ws.Range("H40").Value = 1304.2667
$ws.Range("I40").Value = 1297.4286
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 1297.4286
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -1161.4286
$ws.Range("N40").Value = -1672

$ws.Range("H93").Value = 2398.45
$ws.Range("I93").Value = 1972.2222
$ws.Range("K93").Value = 1972.2222
$ws.Range("M93").Value = -724.2221999999999

$ws.Range("H126").Value = 2358.0908
$ws.Range("I126").Value = 2351.4443
$ws.Range("K126").Value = 7054.3329
$ws.Range("M126").Value = -4584.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 50175.8
$ws.Range("I51").Value = 11069
$ws.Range("J51").Value = 59952.5
$ws.Range("K51").Value = 11069
$ws.Range("L51").Value = 59952.5
$ws.Range("M51").Value = -10559
$ws.Range("N51").Value = -60972.5

$ws.Range("H103").Value = 58994.5
$ws.Range("J103").Value = 58994.5
$ws.Range("L103").Value = 58994.5
$ws.Range("N103").Value = -61338.5

$ws.Range("H104").Value = 91899
$ws.Range("J104").Value = 91899
$ws.Range("L104").Value = 91899
$ws.Range("N104").Value = -98887

$ws.Range("H107").Value = 763.26666
$ws.Range("I107").Value = 738.1111
$ws.Range("J107").Value = 801
$ws.Range("K107").Value = 2214.3333
$ws.Range("L107").Value = 2403
$ws.Range("M107").Value = -294.3332999999998
$ws.Range("N107").Value = -6243

$ws.Range("H132").Value = 2199.3674
$ws.Range("I132").Value = 1908.0217
$ws.Range("K132").Value = 5724.0651
$ws.Range("M132").Value = -3194.0651
